$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing record ("Test Tse" -> "Mimo crimo") ---
$ws.Range("A2").Value = "Mimo crimo"
$ws.Range("B2").Value = "ad646456"
$ws.Range("C2").Value = "'365165415640064545465654"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "a46d"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "797/DR TADLA"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 6000
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 5400

# --- Row 3: update existing record ("Ahmed laaraichi" -> "Ali Ali") ---
$ws.Range("A3").Value = "Ali Ali"
$ws.Range("B3").Value = "aa654556"
$ws.Range("C3").Value = "'566564505100516654656545"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "ad66"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "010/DR010/AV"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2000

# --- Row 4: used to be the totals row, now becomes a normal data record ---
$ws.Range("A4").Value = "Mohamed berrada"
$ws.Range("B4").Value = "ada666"
$ws.Range("C4").Value = "'032165467887646545640545"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "aa121"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "001/DR TANGER/AV"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 4500

# --- Row 5: new record ("IBM"), contract re-activated after suspension ---
$ws.Range("A5").Value = "IBM"
$ws.Range("B5").Value = "'117946464"
$ws.Range("C5").Value = "'114879877777777777777777"
$ws.Range("D5").Value = "BMCE TESTT"
$ws.Range("E5").Value = "BMCE"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "000/DR DEV"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 60000
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 58500

# --- Row 6: new totals row, sums of rows 2-5 ---
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 73000
$ws.Range("J6").Value = 2600
$ws.Range("K6").Value = 70400
